$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 124
$ws.Range("I9").Value = 163.57143
$ws.Range("J9").Value = 54.75
$ws.Range("K9").Value = 163.57143
$ws.Range("L9").Value = 54.75
$ws.Range("M9").Value = 5.428570000000008
$ws.Range("N9").Value = -392.75
$ws.Range("H28").Value = 675.1053000000001
$ws.Range("I28").Value = 484.58334
$ws.Range("J28").Value = 1001.7143
$ws.Range("K28").Value = 484.58334
$ws.Range("L28").Value = 1001.7143
$ws.Range("M28").Value = 0.4166599999999789
$ws.Range("N28").Value = -1971.7143
$ws.Range("H33").Value = 739.11536
$ws.Range("I33").Value = 765.88
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 765.88
$ws.Range("L33").Value = 70
$ws.Range("M33").Value = -536.88
$ws.Range("N33").Value = -528
$ws.Range("H40").Value = 2317.6
$ws.Range("J40").Value = 3394.8572
$ws.Range("L40").Value = 3394.8572
$ws.Range("N40").Value = -3744.8572
$ws.Range("H41").Value = 797.2
$ws.Range("I41").Value = 480
$ws.Range("J41").Value = 832.44446
$ws.Range("K41").Value = 480
$ws.Range("L41").Value = 832.44446
$ws.Range("M41").Value = -40
$ws.Range("N41").Value = -1712.44446
$ws.Range("H43").Value = 989.4
$ws.Range("I43").Value = 650
$ws.Range("J43").Value = 1074.25
$ws.Range("K43").Value = 650
$ws.Range("L43").Value = 1074.25
$ws.Range("M43").Value = -581
$ws.Range("N43").Value = -1212.25
$ws.Range("H53").Value = 586.41174
$ws.Range("I53").Value = 243.25
$ws.Range("J53").Value = 1410
$ws.Range("K53").Value = 243.25
$ws.Range("L53").Value = 1410
$ws.Range("M53").Value = 393.75
$ws.Range("N53").Value = -2684
$ws.Range("H55").Value = 375.08334
$ws.Range("I55").Value = 364.8
$ws.Range("J55").Value = 392.22223
$ws.Range("K55").Value = 364.8
$ws.Range("L55").Value = 392.22223
$ws.Range("M55").Value = -150.8
$ws.Range("N55").Value = -820.2222300000001
$ws.Range("H74").Value = 3951.8462
$ws.Range("I74").Value = 3798
$ws.Range("J74").Value = 4464.6665
$ws.Range("K74").Value = 3798
$ws.Range("L74").Value = 4464.6665
$ws.Range("M74").Value = -2862
$ws.Range("N74").Value = -6336.6665
$ws.Range("H77").Value = 3951.8462
$ws.Range("I77").Value = 3798
$ws.Range("J77").Value = 4464.6665
$ws.Range("K77").Value = 18990
$ws.Range("L77").Value = 22323.3325
$ws.Range("M77").Value = -14310
$ws.Range("N77").Value = -31683.3325
$ws.Range("H127").Value = 1051.9333
$ws.Range("I127").Value = 1023.3333
$ws.Range("J127").Value = 1166.3334
$ws.Range("K127").Value = 3069.9999
$ws.Range("L127").Value = 3499.0002
$ws.Range("M127").Value = 1890.0001
$ws.Range("N127").Value = -13419.0002
$ws.Range("H132").Value = 7459.8374
$ws.Range("I132").Value = 6562.7666
$ws.Range("J132").Value = 9530
$ws.Range("K132").Value = 19688.2998
$ws.Range("L132").Value = 28590
$ws.Range("M132").Value = -17158.2998
$ws.Range("N132").Value = -33650
$ws.Range("H138").Value = 1411.909
$ws.Range("I138").Value = 1078.6897
$ws.Range("J138").Value = 2429.1052
$ws.Range("K138").Value = 3236.0691
$ws.Range("L138").Value = 7287.3156
$ws.Range("M138").Value = 1903.9309
$ws.Range("N138").Value = -17567.3156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6641.65
$ws.Range("I32").Value = 5640.053
$ws.Range("J32").Value = 22333.334
$ws.Range("K32").Value = 5640.053
$ws.Range("L32").Value = 22333.334
$ws.Range("M32").Value = -5353.053
$ws.Range("N32").Value = -22907.334
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1349.6666
$ws.Range("I20").Value = 959.5625
$ws.Range("J20").Value = 1716.8235
$ws.Range("K20").Value = 959.5625
$ws.Range("L20").Value = 1716.8235
$ws.Range("M20").Value = -712.5625
$ws.Range("N20").Value = -2210.8235
$ws.Range("H94").Value = 2272.3333
$ws.Range("I94").Value = 2406.182
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 2406.182
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -1955.182
$ws.Range("N94").Value = -1702
$ws.Range("H134").Value = 3067.7021
$ws.Range("I134").Value = 1916.3889
$ws.Range("J134").Value = 3782.3103
$ws.Range("K134").Value = 5749.1667
$ws.Range("L134").Value = 11346.9309
$ws.Range("M134").Value = -3214.1667
$ws.Range("N134").Value = -16416.9309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4330967.5
$ws.Range("I31").Value = 1227.8793
$ws.Range("J31").Value = 17548068
$ws.Range("K31").Value = 1227.8793
$ws.Range("L31").Value = 17548068
$ws.Range("M31").Value = -932.8793000000001
$ws.Range("N31").Value = -17548658
$ws.Range("H34").Value = 4330967.5
$ws.Range("I34").Value = 1227.8793
$ws.Range("J34").Value = 17548068
$ws.Range("K34").Value = 1227.8793
$ws.Range("L34").Value = 17548068
$ws.Range("M34").Value = -1025.8793
$ws.Range("N34").Value = -17548472
$ws.Range("H43").Value = 21114.25
$ws.Range("J43").Value = 21114.25
$ws.Range("L43").Value = 21114.25
$ws.Range("N43").Value = -21482.25
$ws.Range("H101").Value = 21114.25
$ws.Range("J101").Value = 21114.25
$ws.Range("L101").Value = 21114.25
$ws.Range("N101").Value = -27604.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 271.9091
$ws.Range("I2").Value = 21.9
$ws.Range("J2").Value = 480.25
$ws.Range("K2").Value = 131.4
$ws.Range("L2").Value = 2881.5
$ws.Range("M2").Value = -18.39999999999998
$ws.Range("N2").Value = -3107.5
$ws.Range("H22").Value = 2406.3809
$ws.Range("I22").Value = 2283.5
$ws.Range("J22").Value = 2455.5334
$ws.Range("K22").Value = 6850.5
$ws.Range("L22").Value = 7366.600199999999
$ws.Range("M22").Value = -6681.5
$ws.Range("N22").Value = -7704.600199999999
$ws.Range("H27").Value = 2406.3809
$ws.Range("I27").Value = 2283.5
$ws.Range("J27").Value = 2455.5334
$ws.Range("K27").Value = 6850.5
$ws.Range("L27").Value = 7366.600199999999
$ws.Range("M27").Value = -6748.5
$ws.Range("N27").Value = -7570.600199999999
$ws.Range("H50").Value = 86.27273
$ws.Range("I50").Value = 44.9
$ws.Range("J50").Value = 500
$ws.Range("K50").Value = 134.7
$ws.Range("L50").Value = 1500
$ws.Range("M50").Value = 346.3
$ws.Range("N50").Value = -2462
$ws.Range("H53").Value = 86.27273
$ws.Range("I53").Value = 44.9
$ws.Range("J53").Value = 500
$ws.Range("K53").Value = 134.7
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = 346.3
$ws.Range("N53").Value = -2462
$ws.Range("H80").Value = 3298.2856
$ws.Range("J80").Value = 3298.2856
$ws.Range("L80").Value = 9894.856800000001
$ws.Range("N80").Value = -11766.8568
$ws.Range("H83").Value = 3298.2856
$ws.Range("J83").Value = 3298.2856
$ws.Range("L83").Value = 29684.5704
$ws.Range("N83").Value = -39044.5704
$ws.Range("H86").Value = 1369.6
$ws.Range("I86").Value = 616
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1848
$ws.Range("L86").Value = 7500
$ws.Range("M86").Value = -662
$ws.Range("N86").Value = -9872
$ws.Range("H89").Value = 1369.6
$ws.Range("I89").Value = 616
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 5544
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = 384
$ws.Range("N89").Value = -34356
$ws.Range("H131").Value = 870
$ws.Range("I131").Value = 440.76923
$ws.Range("J131").Value = 1218.75
$ws.Range("K131").Value = 1322.30769
$ws.Range("L131").Value = 3656.25
$ws.Range("M131").Value = 3717.69231
$ws.Range("N131").Value = -13736.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 962
$ws.Range("I97").Value = 936.6667
$ws.Range("K97").Value = 936.6667
$ws.Range("M97").Value = -440.6667
$ws.Range("H132").Value = 3355.5483
$ws.Range("I132").Value = 3121
$ws.Range("J132").Value = 3484.55
$ws.Range("K132").Value = 9363
$ws.Range("L132").Value = 10453.65
$ws.Range("M132").Value = -6833
$ws.Range("N132").Value = -15513.65

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 221.66667
$ws.Range("J55").Value = 333.33334
$ws.Range("L55").Value = 333.33334
$ws.Range("N55").Value = -679.33334
$ws.Range("H82").Value = 1543.8148
$ws.Range("I82").Value = 1032.3889
$ws.Range("J82").Value = 2566.6667
$ws.Range("K82").Value = 1032.3889
$ws.Range("L82").Value = 2566.6667
$ws.Range("M82").Value = -671.3888999999999
$ws.Range("N82").Value = -3288.6667
$ws.Range("H85").Value = 1543.8148
$ws.Range("I85").Value = 1032.3889
$ws.Range("J85").Value = 2566.6667
$ws.Range("K85").Value = 1032.3889
$ws.Range("L85").Value = 2566.6667
$ws.Range("M85").Value = 215.6111000000001
$ws.Range("N85").Value = -5062.6667
$ws.Range("H98").Value = 26999.5
$ws.Range("J98").Value = 26999.5
$ws.Range("L98").Value = 26999.5
$ws.Range("N98").Value = -32989.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1550.4615
$ws.Range("I81").Value = 1513
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3026
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -1965
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1550.4615
$ws.Range("I84").Value = 1513
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 15130
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -9826
$ws.Range("N84").Value = -30608
